$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "63.487.46"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  +3.44%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.065.38"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  +2.08%  "

$ws.Range("E4").Value = "  +0.02%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "549.51"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +2.09%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "140.33"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +3.21%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +0.03%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.055.39"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +1.83%  "

$ws.Range("E9").Value = "  +1.14%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.49"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +5.95%  "

$ws.Range("E11").Value = "  +1.00%  "

$ws.Range("E12").Value = "  +1.30%  "

$ws.Range("E13").Value = "  +2.46%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "34.82"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +1.50%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.563.51"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +2.12%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "63.436.81"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +3.28%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.068.60"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +2.14%  "

$ws.Range("E18").Value = "  -1.25%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.74"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +1.61%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "482.62"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +2.91%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.69"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +3.07%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.671"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -0.91%  "

$ws.Range("E23").Value = "  +4.32%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "80.80"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +1.28%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "12.63"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +4.68%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.998"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -0.08%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.75"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +2.18%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.91"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -0.68%  "

$ws.Range("E29").Value = "  +5.03%  "

$ws.Range("E30").Value = "  -0.09%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "26.12"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +1.81%  "

$ws.Range("E32").Value = "  -0.13%  "

$ws.Range("E33").Value = "  +7.32%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.71"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +3.70%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "55.48"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -0.56%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.98"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +1.33%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "468.19"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +2.84%  "

$ws.Range("E38").Value = "  +3.86%  "

$ws.Range("E39").Value = "  +3.04%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.068.17"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -4.39%  "

$ws.Range("E41").Value = "  +0.50%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.25"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +1.04%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.57"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +2.78%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "27.93"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +0.04%  "

$ws.Range("E45").Value = "  +3.38%  "

$ws.Range("E47").Value = "  +2.25%  "

$ws.Range("E48").Value = "  +1.32%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "116.42"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -3.19%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0₃0509"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +1.89%  "

$ws.Range("E51").Value = "  +3.00%  "
